$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the existing column filter (SIVhu) and unhide all filtered rows,
# matching the removal of <filterColumn>/<filters> and row hidden="1" attrs
# while keeping the sortState in the autoFilter definition.
$ws.ShowAllData()

# Append rows for the new August 2020 sequencing run samples (rows 176-195)
$ws.Range("A176").Value = 'HIV_2_Group_F_1st_Generation_Mouse_1_W11'
$ws.Range("B176").Value = 'HIV'
$ws.Range("C176").Value = 1
$ws.Range("D176").Value = 11
$ws.Range("E176").Value = 1
$ws.Range("F176").Value = 'A'

$ws.Range("A177").Value = 'HIV_2_Group_F_1st_Generation_Mouse_1_W3'
$ws.Range("B177").Value = 'HIV'
$ws.Range("C177").Value = 1
$ws.Range("D177").Value = 3
$ws.Range("E177").Value = 1
$ws.Range("F177").Value = 'A'

$ws.Range("A178").Value = 'SIVB670_2nd_Generation_Mouse_J2873_W23'
$ws.Range("B178").Value = 'SIVB670'
$ws.Range("C178").Value = 2
$ws.Range("D178").Value = 23
$ws.Range("E178").Value = 'J2873'
$ws.Range("F178").Value = 'A'

$ws.Range("A179").Value = 'SIVB670_2nd_Generation_Mouse_J2878_W23'
$ws.Range("B179").Value = 'SIVB670'
$ws.Range("C179").Value = 2
$ws.Range("D179").Value = 23
$ws.Range("E179").Value = 'J2878'
$ws.Range("F179").Value = 'B'

$ws.Range("A180").Value = 'SIVB670_3rd_Generation_Mouse_B1076_W11'
$ws.Range("B180").Value = 'SIVB670'
$ws.Range("C180").Value = 3
$ws.Range("D180").Value = 11
$ws.Range("E180").Value = 'B1076'
$ws.Range("F180").Value = 'A'

$ws.Range("A181").Value = 'SIVB670_3rd_Generation_Mouse_B1076_W23'
$ws.Range("B181").Value = 'SIVB670'
$ws.Range("C181").Value = 3
$ws.Range("D181").Value = 23
$ws.Range("E181").Value = 'B1076'
$ws.Range("F181").Value = 'A'

$ws.Range("A182").Value = 'SIVB670_3rd_Generation_Mouse_B1076_W3'
$ws.Range("B182").Value = 'SIVB670'
$ws.Range("C182").Value = 3
$ws.Range("D182").Value = 3
$ws.Range("E182").Value = 'B1076'
$ws.Range("F182").Value = 'A'

$ws.Range("A183").Value = 'SIVB670_3rd_Generation_Mouse_B1077_W23'
$ws.Range("B183").Value = 'SIVB670'
$ws.Range("C183").Value = 3
$ws.Range("D183").Value = 23
$ws.Range("E183").Value = 'B1077'
$ws.Range("F183").Value = 'B'

$ws.Range("A184").Value = 'SIVB670_3rd_Generation_Mouse_B1077_W3'
$ws.Range("B184").Value = 'SIVB670'
$ws.Range("C184").Value = 3
$ws.Range("D184").Value = 3
$ws.Range("E184").Value = 'B1077'
$ws.Range("F184").Value = 'B'

$ws.Range("A185").Value = 'SIVcpzEK505_3rd_Generation_Mouse_2404_W3'
$ws.Range("B185").Value = 'SIVcpzEK505'
$ws.Range("C185").Value = 3
$ws.Range("D185").Value = 3
$ws.Range("E185").Value = 'J2404'
$ws.Range("F185").Value = 'A'

$ws.Range("A186").Value = 'SIVcpzLB715_3rd_Generation_Mouse_J12_W11'
$ws.Range("B186").Value = 'SIVcpzLB715'
$ws.Range("C186").Value = 3
$ws.Range("D186").Value = 11
$ws.Range("E186").Value = 'J12'
$ws.Range("F186").Value = 'A'

$ws.Range("A187").Value = 'SIVcpzLB715_3rd_Generation_Mouse_J12_W19'
$ws.Range("B187").Value = 'SIVcpzLB715'
$ws.Range("C187").Value = 3
$ws.Range("D187").Value = 19
$ws.Range("E187").Value = 'J12'
$ws.Range("F187").Value = 'A'

$ws.Range("A188").Value = 'SIVcpzLB715_3rd_Generation_Mouse_J12_W22'
$ws.Range("B188").Value = 'SIVcpzLB715'
$ws.Range("C188").Value = 3
$ws.Range("D188").Value = 22
$ws.Range("E188").Value = 'J12'
$ws.Range("F188").Value = 'A'

$ws.Range("A189").Value = 'SIVcpzLB715_3rd_Generation_Mouse_J12_W3'
$ws.Range("B189").Value = 'SIVcpzLB715'
$ws.Range("C189").Value = 3
$ws.Range("D189").Value = 3
$ws.Range("E189").Value = 'J12'
$ws.Range("F189").Value = 'A'

$ws.Range("A190").Value = 'SIVcpzLB715_3rd_Generation_Mouse_J4_W11'
$ws.Range("B190").Value = 'SIVcpzLB715'
$ws.Range("C190").Value = 3
$ws.Range("D190").Value = 11
$ws.Range("E190").Value = 'J4'
$ws.Range("F190").Value = 'B'

$ws.Range("A191").Value = 'SIVcpzLB715_3rd_Generation_Mouse_J4_W19'
$ws.Range("B191").Value = 'SIVcpzLB715'
$ws.Range("C191").Value = 3
$ws.Range("D191").Value = 19
$ws.Range("E191").Value = 'J4'
$ws.Range("F191").Value = 'B'

$ws.Range("A192").Value = 'SIVcpzLB715_3rd_Generation_Mouse_J4_W3'
$ws.Range("B192").Value = 'SIVcpzLB715'
$ws.Range("C192").Value = 3
$ws.Range("D192").Value = 3
$ws.Range("E192").Value = 'J4'
$ws.Range("F192").Value = 'B'

$ws.Range("A193").Value = 'SIVcpzMB897_3rd_Generation_Mouse_J2904_W11'
$ws.Range("B193").Value = 'SIVcpzMB897'
$ws.Range("C193").Value = 3
$ws.Range("D193").Value = 11
$ws.Range("E193").Value = 'J2904'
$ws.Range("F193").Value = 'A'

$ws.Range("A194").Value = 'SIVcpzMB897_3rd_Generation_Mouse_J2904_W18'
$ws.Range("B194").Value = 'SIVcpzMB897'
$ws.Range("C194").Value = 3
$ws.Range("D194").Value = 18
$ws.Range("E194").Value = 'J2904'
$ws.Range("F194").Value = 'A'

$ws.Range("A195").Value = 'SIVcpzMB897_3rd_Generation_Mouse_J2904_W3'
$ws.Range("B195").Value = 'SIVcpzMB897'
$ws.Range("C195").Value = 3
$ws.Range("D195").Value = 3
$ws.Range("E195").Value = 'J2904'
$ws.Range("F195").Value = 'A'

# Scroll/selection to the newly appended rows, matching the updated sheetView
[void]$ws.Activate()
$excel.ActiveWindow.ScrollRow = 167
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("A176").Select()
